$d = $word.ActiveDocument

# Locate the paragraph ending with "...end to end." (the "3 years of
# production-grade..." sentence) so the new paragraph is anchored by
# content rather than a brittle hard-coded index.
$anchor = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "I have 3 years of production-grade*") {
        $anchor = $para
        break
    }
}

# Add a brand new paragraph right after it.
$anchor.Range.InsertParagraphAfter()
$newPara = $anchor.Next()

# First run of the new paragraph: the sentence without its closing period.
$newPara.Range.Text = "I have worked in a fast paced startup and more organized MNC"

# Second run: just the period. Inserting it via a collapsed Range right
# after the last real character (End - 1, since End itself sits on the
# boundary shared with the following paragraph) and nudging a character
# property on/off forces it to land in its own run instead of being
# merged back into the first one.
$insertAt = $newPara.Range.Start + $newPara.Range.Text.Length - 1
$periodRange = $d.Range($insertAt, $insertAt)
$periodRange.InsertAfter(".")
$periodRange.Bold = 1
$periodRange.Bold = 0
